$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be introduced in this exact order so that the
# --- resulting sharedStrings table assigns indices 70,71,72,73 as expected:
#   70 "not full string, just prot change part"
#   71 "could be bug, mix up btwn variant reads, sample reads"
#   72 "exac"
#   73 "failing hard"

# Row 43: cosmic id row gains "yes" / "cosmic" columns, and the comment in E
# changes from the generic cosmic-db-needs-updating note to this more
# specific one (first brand-new string -> index 70).
$ws.Range("C43").Value = "yes"
$ws.Range("D43").Value = "cosmic"
$ws.Range("E43").Value = "not full string, just prot change part"

# Row 42 ("target") also gets the "cosmic" comment in column D.
$ws.Range("D42").Value = "cosmic"

# Rows 16 & 17 get a new comment in column E (second brand-new string ->
# index 71).
$ws.Range("E16").Value = "could be bug, mix up btwn variant reads, sample reads"
$ws.Range("E17").Value = "could be bug, mix up btwn variant reads, sample reads"

# Column A values added/updated for a few test-case grouping numbers.
$ws.Range("A20").Value = 3
$ws.Range("A21").Value = 3
$ws.Range("A37").Value = 2
$ws.Range("A38").Value = 2

# New ExAC test cases (rows 48-56): A becomes 1 (was 2), plus new columns
# C/D/E ("yes" / "exac" / "failing hard" - third & fourth brand-new strings
# -> indices 72 and 73).
$exacRows = 48..56
foreach ($r in $exacRows) {
    $ws.Range("A$r").Value = 1
    $ws.Range("C$r").Value = "yes"
    $ws.Range("D$r").Value = "exac"
    $ws.Range("E$r").Value = "failing hard"
}

# --- Sheet view / window state ---
$excel.ActiveWindow.Zoom = 130
try {
    $excel.ActiveWindow.ScrollRow = 31
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Older/limited hosts may not support window scroll positioning; not
    # critical to the data change, so ignore.
}
$ws.Range("E48:E56").Select()
